$wb = $excel.ActiveWorkbook

# --- Move the current selection on AlertsTemplateManagement to G10 and
#     deactivate it, so the new sheet can become the active tab. ---
$atm = $wb.Worksheets.Item("AlertsTemplateManagement")
[void]$atm.Activate()
[void]$atm.Range("G10").Select()

# --- Insert a new worksheet "AlertsNotificationManagement" right after
#     "AlertsTemplateManagement" (becomes sheetId 16 / rId5, all later
#     sheet rIds shift by one automatically). ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $atm)
$newSheet.Name = "AlertsNotificationManagement"

# --- Bring over the header/value cell formatting (styles 54/53/55) from
#     the structurally-identical CoreAlertsPlaceholderManagement sheet. ---
$src1 = $wb.Worksheets.Item("CoreAlertsPlaceholderManagement")
$src1.Range("A1:B2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# --- Bring over the last-column header/value formatting (styles 57/58)
#     from AlertsTemplateManagement's own C column. ---
$src2 = $wb.Worksheets.Item("AlertsTemplateManagement")
$src2.Range("C1:C2").Copy()
$newSheet.Range("C1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Populate the new sheet's data. ---
$newSheet.Range("A1").Value = "TestScenario"
$newSheet.Range("B1").Value = "Run"
$newSheet.Range("C1").Value = "ScheduleType"
$newSheet.Range("A2").Value = "AlertsNotificationManagement"
$newSheet.Range("B2").Value = "Yes"
$newSheet.Range("C2").Value = "One Time"

# --- Column widths matching the authored sheet. ---
$newSheet.Columns("A").ColumnWidth = 26.73
$newSheet.Columns("C").ColumnWidth = 17.82

# --- New sheet becomes the active tab with its own selection. ---
[void]$newSheet.Activate()
[void]$newSheet.Range("C9").Select()

Write-Output "AlertsNotificationManagement sheet added."
